$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 5.824850660931335
$ws.Cells.Item(2,4).Value = 6.074095176849792
$ws.Cells.Item(2,5).Value = 12.40051032086643
$ws.Cells.Item(2,6).Value = 48.03090014599303
$ws.Cells.Item(2,7).Value = 3.731625644699826
$ws.Cells.Item(2,9).Value = 33.93957203742939
$ws.Cells.Item(2,10).Value = 10.47764372969178
$ws.Cells.Item(2,11).Value = 21.64679668001352
$ws.Cells.Item(3,3).Value = 5.839047525712984
$ws.Cells.Item(3,4).Value = 6.0575534092765
$ws.Cells.Item(3,5).Value = 12.36843643059752
$ws.Cells.Item(3,6).Value = 47.93721475416519
$ws.Cells.Item(3,7).Value = 3.736002254940679
$ws.Cells.Item(3,9).Value = 33.90180973255217
$ws.Cells.Item(3,10).Value = 10.4893810144634
$ws.Cells.Item(3,11).Value = 21.26327532129198
$ws.Cells.Item(4,3).Value = 5.848095460125783
$ws.Cells.Item(4,4).Value = 6.047926270026275
$ws.Cells.Item(4,5).Value = 12.35130689126168
$ws.Cells.Item(4,6).Value = 47.89442809230174
$ws.Cells.Item(4,7).Value = 3.738824997905463
$ws.Cells.Item(4,9).Value = 33.88886290587809
$ws.Cells.Item(4,10).Value = 10.49852316796066
$ws.Cells.Item(4,11).Value = 21.03009222319272
$ws.Cells.Item(5,3).Value = 5.851866447008917
$ws.Cells.Item(5,4).Value = 6.044139194146917
$ws.Cells.Item(5,5).Value = 12.34497530149623
$ws.Cells.Item(5,6).Value = 47.88069105763327
$ws.Cells.Item(5,7).Value = 3.7400095004851
$ws.Cells.Item(5,9).Value = 33.88615270217656
$ws.Cells.Item(5,10).Value = 10.50273439384826
$ws.Cells.Item(5,11).Value = 20.93577792114452
$ws.Cells.Item(6,3).Value = 5.85249770168625
$ws.Cells.Item(6,4).Value = 6.043518654295088
$ws.Cells.Item(6,5).Value = 12.34396324234868
$ws.Cells.Item(6,6).Value = 47.87863315820938
$ws.Cells.Item(6,7).Value = 3.740208256734127
$ws.Cells.Item(6,9).Value = 33.88585732725387
$ws.Cells.Item(6,10).Value = 10.50346297121026
$ws.Cells.Item(6,11).Value = 20.92016380457008
$ws.Cells.Item(7,3).Value = 5.848145976482455
$ws.Cells.Item(7,4).Value = 6.047874641505116
$ws.Cells.Item(7,5).Value = 12.35121886910915
$ws.Cells.Item(7,6).Value = 47.89422786510499
$ws.Cells.Item(7,7).Value = 3.738840833807969
$ws.Cells.Item(7,9).Value = 33.88881598022472
$ws.Cells.Item(7,10).Value = 10.49857799667244
$ws.Cells.Item(7,11).Value = 21.02881722044231
$ws.Cells.Item(8,3).Value = 5.829677429572637
$ws.Cells.Item(8,4).Value = 6.068283416174697
$ws.Cells.Item(8,5).Value = 12.3889217322514
$ws.Cells.Item(8,6).Value = 47.9955360236205
$ws.Cells.Item(8,7).Value = 3.733106665819584
$ws.Cells.Item(8,9).Value = 33.92442169039558
$ws.Cells.Item(8,10).Value = 10.48128850937104
$ws.Cells.Item(8,11).Value = 21.514159604631
$ws.Cells.Item(9,3).Value = 5.796057666798779
$ws.Cells.Item(9,4).Value = 6.112378023120883
$ws.Cells.Item(9,5).Value = 12.48299552346637
$ws.Cells.Item(9,6).Value = 48.31124546402459
$ws.Cells.Item(9,7).Value = 3.722930483941101
$ws.Cells.Item(9,9).Value = 34.07578517364097
$ws.Cells.Item(9,10).Value = 10.46277959495353
$ws.Cells.Item(9,11).Value = 22.47845036753195
$ws.Cells.Item(10,3).Value = 5.772899557118504
$ws.Cells.Item(10,4).Value = 6.147082441097069
$ws.Cells.Item(10,5).Value = 12.5640784427676
$ws.Cells.Item(10,6).Value = 48.61444520773555
$ws.Cells.Item(10,7).Value = 3.716096286895096
$ws.Cells.Item(10,9).Value = 34.23694853599027
$ws.Cells.Item(10,10).Value = 10.45861507979384
$ws.Cells.Item(10,11).Value = 23.18705276476855
$ws.Cells.Item(11,3).Value = 5.762690537485754
$ws.Cells.Item(11,4).Value = 6.163333329069923
$ws.Cells.Item(11,5).Value = 12.60348599037301
$ws.Cells.Item(11,6).Value = 48.76773869789787
$ws.Cells.Item(11,7).Value = 3.71312472908655
$ws.Cells.Item(11,9).Value = 34.32111313541599
$ws.Cells.Item(11,10).Value = 10.45877693659414
$ws.Cells.Item(11,11).Value = 23.50792395778221
$ws.Cells.Item(12,3).Value = 5.758870813279369
$ws.Cells.Item(12,4).Value = 6.169550473509925
$ws.Cells.Item(12,5).Value = 12.61876388246456
$ws.Cells.Item(12,6).Value = 48.82798025290928
$ws.Cells.Item(12,7).Value = 3.712019076801348
$ws.Cells.Item(12,9).Value = 34.35453988607131
$ws.Cells.Item(12,10).Value = 10.45913430943497
$ws.Cells.Item(12,11).Value = 23.62909697917176
$ws.Cells.Item(13,3).Value = 5.759691413979666
$ws.Cells.Item(13,4).Value = 6.168208739914682
$ws.Cells.Item(13,5).Value = 12.61545784823889
$ws.Cells.Item(13,6).Value = 48.8149089529546
$ws.Cells.Item(13,7).Value = 3.71225632875781
$ws.Cells.Item(13,9).Value = 34.34727176327189
$ws.Cells.Item(13,10).Value = 10.45904416983584
$ws.Cells.Item(13,11).Value = 23.60301693052246
$ws.Cells.Item(14,3).Value = 5.762375363946768
$ws.Cells.Item(14,4).Value = 6.163843568483684
$ws.Cells.Item(14,5).Value = 12.60473584029988
$ws.Cells.Item(14,6).Value = 48.77265097539746
$ws.Cells.Item(14,7).Value = 3.713033374146864
$ws.Cells.Item(14,9).Value = 34.32383203794422
$ws.Cells.Item(14,10).Value = 10.45880040224637
$ws.Cells.Item(14,11).Value = 23.5179003048381
$ws.Cells.Item(15,3).Value = 5.764025360117452
$ws.Cells.Item(15,4).Value = 6.161177915414519
$ws.Cells.Item(15,5).Value = 12.59821431160653
$ws.Cells.Item(15,6).Value = 48.74705171406239
$ws.Cells.Item(15,7).Value = 3.713511886843476
$ws.Cells.Item(15,9).Value = 34.3096769056159
$ws.Cells.Item(15,10).Value = 10.45868965516057
$ws.Cells.Item(15,11).Value = 23.46571690928229
$ws.Cells.Item(16,3).Value = 5.773573241089763
$ws.Cells.Item(16,4).Value = 6.146029448714771
$ws.Cells.Item(16,5).Value = 12.56155318609782
$ws.Cells.Item(16,6).Value = 48.60473496758548
$ws.Cells.Item(16,7).Value = 3.716293237190019
$ws.Cells.Item(16,9).Value = 34.23166616702312
$ws.Cells.Item(16,10).Value = 10.4586459164898
$ws.Cells.Item(16,11).Value = 23.1660428497029
$ws.Cells.Item(17,3).Value = 5.77951353329108
$ws.Cells.Item(17,4).Value = 6.136852713050656
$ws.Cells.Item(17,5).Value = 12.53970358740651
$ws.Cells.Item(17,6).Value = 48.52135315942053
$ws.Cells.Item(17,7).Value = 3.718034585921711
$ws.Cells.Item(17,9).Value = 34.18658556965005
$ws.Cells.Item(17,10).Value = 10.45914607380202
$ws.Cells.Item(17,11).Value = 22.9817362031124
$ws.Cells.Item(18,3).Value = 5.782960927920538
$ws.Cells.Item(18,4).Value = 6.131618345939846
$ws.Cells.Item(18,5).Value = 12.52737419746161
$ws.Cells.Item(18,6).Value = 48.47484179531671
$ws.Cells.Item(18,7).Value = 3.719049100245628
$ws.Cells.Item(18,9).Value = 34.1616781676568
$ws.Cells.Item(18,10).Value = 10.45962728689047
$ws.Cells.Item(18,11).Value = 22.87559595139098
$ws.Cells.Item(19,3).Value = 5.78413344771089
$ws.Cells.Item(19,4).Value = 6.129853712887806
$ws.Cells.Item(19,5).Value = 12.52324076509871
$ws.Cells.Item(19,6).Value = 48.45934293202622
$ws.Cells.Item(19,7).Value = 3.71939482347063
$ws.Cells.Item(19,9).Value = 34.15342054367257
$ws.Cells.Item(19,10).Value = 10.45982344332568
$ws.Cells.Item(19,11).Value = 22.83963982542237
$ws.Cells.Item(20,3).Value = 5.778878006726749
$ws.Cells.Item(20,4).Value = 6.137825078512553
$ws.Cells.Item(20,5).Value = 12.54200494417301
$ws.Cells.Item(20,6).Value = 48.53007957959692
$ws.Cells.Item(20,7).Value = 3.717847878496103
$ws.Cells.Item(20,9).Value = 34.1912787543334
$ws.Cells.Item(20,10).Value = 10.45907279775493
$ws.Cells.Item(20,11).Value = 23.00137051572551
$ws.Cells.Item(21,3).Value = 5.761585773641021
$ws.Cells.Item(21,4).Value = 6.16512403477737
$ws.Cells.Item(21,5).Value = 12.60787557964511
$ws.Cells.Item(21,6).Value = 48.78500381636373
$ws.Cells.Item(21,7).Value = 3.712804605919082
$ws.Cells.Item(21,9).Value = 34.33067469070225
$ws.Cells.Item(21,10).Value = 10.45886396474856
$ws.Cells.Item(21,11).Value = 23.54291112597103
$ws.Cells.Item(22,3).Value = 5.750553357947067
$ws.Cells.Item(22,4).Value = 6.183332763145319
$ws.Cells.Item(22,5).Value = 12.65299181048327
$ws.Cells.Item(22,6).Value = 48.96438112969809
$ws.Cells.Item(22,7).Value = 3.709622786033388
$ws.Cells.Item(22,9).Value = 34.43084023708508
$ws.Cells.Item(22,10).Value = 10.46045334934112
$ws.Cells.Item(22,11).Value = 23.89484043690094
$ws.Cells.Item(23,3).Value = 5.756417153505931
$ws.Cells.Item(23,4).Value = 6.173581944489282
$ws.Cells.Item(23,5).Value = 12.62872604078108
$ws.Cells.Item(23,6).Value = 48.86748234820718
$ws.Cells.Item(23,7).Value = 3.71131057516683
$ws.Cells.Item(23,9).Value = 34.37655307507089
$ws.Cells.Item(23,10).Value = 10.45944705599999
$ws.Cells.Item(23,11).Value = 23.70723075568954
$ws.Cells.Item(24,3).Value = 5.779165227814099
$ws.Cells.Item(24,4).Value = 6.137385342690449
$ws.Cells.Item(24,5).Value = 12.54096377684514
$ws.Cells.Item(24,6).Value = 48.52612992228211
$ws.Cells.Item(24,7).Value = 3.717932247206317
$ws.Cells.Item(24,9).Value = 34.18915381916985
$ws.Cells.Item(24,10).Value = 10.4591053226602
$ws.Cells.Item(24,11).Value = 22.99249439405969
$ws.Cells.Item(25,3).Value = 5.804878919575198
$ws.Cells.Item(25,4).Value = 6.100031056951727
$ws.Cells.Item(25,5).Value = 12.45541982499531
$ws.Cells.Item(25,6).Value = 48.21330123817023
$ws.Cells.Item(25,7).Value = 3.72556996612258
$ws.Cells.Item(25,9).Value = 34.02607215845769
$ws.Cells.Item(25,10).Value = 10.46613303785358
$ws.Cells.Item(25,11).Value = 22.21702909912886
